# Commit log update + progress report 2
# Appends new commit-log rows (2-9) to the existing Table1 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NAME, DATE, FILE, COMMENT
$rows = @(
    @("Omri", "16.12.18", "address_translation.v",   "interface created"),
    @("Ori",  "16.12.18", "cache.v",                  "interface created"),
    @("Ori",  "16.12.18", "spi.v",                    "deleting unnesecery states from state machine"),
    @("Omri", "16.12.18", "interconnect.v",            "adjesting interface according to address_translation.v"),
    @("Omri", "19.12.18", "per_home_logic.v",          "adjesting interface according to address_translation.v"),
    @("Ori",  "19.12.18", "cpu_if.v",                  "taking down one sample stage for read data"),
    @("Omri", "23.12.18", "address_translation4k.v",   "inserting basic translation for spesific case in order to check functionality"),
    @("Ori",  "23.12.18", "spi.v",                     "inserting NVM read request stage to spi state machine")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# The COMMENT column now holds much longer text -> widen it to fit.
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Leave selection on the last entered cell, like the saved workbook shows.
$ws.Range("D9").Select() | Out-Null
